# Assigned toll classes for major arterials
#
# Adds 20 new rows (274-293) of facility_name / toll class id data to the
# "Inputs_for_tollcalib" sheet, covering several arterials that parallel
# major freeways (US-101, I-880/I-580, I-80, I-280, SR-4, 237/680).
#
# The facility_name values (column A) are written in the exact order that
# first introduces each new shared string, so the resulting shared string
# table indices line up with the canonical workbook. Column B (toll class
# id) is then populated in one shot for rows 274-293.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (facility_name) -------------------------------------------
# Order below = order new distinct strings are first used (NOT row order),
# which preserves the intended shared-string ordering in the saved file.
$ws.Range("A276").Value = '[Santa Clara: paralleling US-101] Monterey Rd - NB'
$ws.Range("A277").Value = '[Santa Clara: paralleling US-101] Monterey Rd - SB'
$ws.Range("A274").Value = '[Alameda: paralleling I-880/I-580] International Blvd (to E 14th St to Mission Blvd) - NB'
$ws.Range("A275").Value = '[Alameda: paralleling I-880/I-580] International Blvd (to E 14th St to Mission Blvd) - SB'
$ws.Range("A278").Value = '[San Francisco: paralleling US-101] 3rd (to Bayshore Blvd to Airport Blvd) - NB'
$ws.Range("A279").Value = '[San Francisco: paralleling US-101] 3rd (to Bayshore Blvd to Airport Blvd) - SB'
$ws.Range("A280").Value = '[Alameda: paralleling I-80] San Pablo - NB'
$ws.Range("A281").Value = '[Alameda: paralleling I-80] San Pablo - SB'
$ws.Range("A292").Value = '[Santa Clara: paralleling I-280] Foothill Expy - NB'
$ws.Range("A293").Value = '[Santa Clara: paralleling I-280] Foothill Expy - SB'
$ws.Range("A290").Value = '[Contra Costa: paralleling SR-4] Leland Rd (to Delta Fair Blvd) - EB'
$ws.Range("A291").Value = '[Contra Costa: paralleling SR-4] Leland Rd (to Delta Fair Blvd) -WB'
$ws.Range("A288").Value = '[Solano: paralleling I-80] Texas St - EB'
$ws.Range("A289").Value = '[Solano: paralleling I-80] Texas St - WB'
$ws.Range("A286").Value = '[Santa Clara: paralleling 237/680] Tasman Dr (to N Capitol Ave to E Capitol Expy) - EB'
$ws.Range("A287").Value = '[Santa Clara: paralleling 237/680] Tasman Dr (to N Capitol Ave to E Capitol Expy) - EB'
$ws.Range("A284").Value = '[Sonoma: paralleling US-101] Mendocino Ave +  Santa Rosa Ave - NB'
$ws.Range("A285").Value = '[Sonoma: paralleling US-101] Mendocino Ave +  Santa Rosa Ave - SB'
$ws.Range("A282").Value = '[San Mateo: paralleling US-101/I-280] El Camino Real - NB'
$ws.Range("A283").Value = '[San Mateo: paralleling US-101/I-280] El Camino Real - SB'

# --- Column B (toll class id) -------------------------------------------
$bVals = New-Object 'object[,]' 20,1
$bVals[0,0]  = 770580
$bVals[1,0]  = 770851
$bVals[2,0]  = 770101
$bVals[3,0]  = 770102
$bVals[4,0]  = 770103
$bVals[5,0]  = 770104
$bVals[6,0]  = 770080
$bVals[7,0]  = 770081
$bVals[8,0]  = 770105
$bVals[9,0]  = 770106
$bVals[10,0] = 770107
$bVals[11,0] = 770108
$bVals[12,0] = 770237
$bVals[13,0] = 770238
$bVals[14,0] = 770082
$bVals[15,0] = 770083
$bVals[16,0] = 770004
$bVals[17,0] = 770005
$bVals[18,0] = 770280
$bVals[19,0] = 770281
$ws.Range("B274:B293").Value = $bVals

# --- Column A width (widen to fit the new, longer facility names) -------
$ws.Columns.Item(1).ColumnWidth = 79

# --- View / selection state ----------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 229
$win.ScrollColumn = 1
$ws.Range("C279").Select()
